$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($row1, $row2, $col) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

# Rows 8 and 9: swap Id, Taxonsorteringsordning, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Ost, Nord, Starttid, Sluttid
$cols89 = @("A","B","E","F","G","H","Q","R","Z","AB")
foreach ($col in $cols89) {
    Swap-Cell 8 9 $col
}

# Rows 13 and 14: swap Id, Ost, Nord, Starttid, Sluttid
$cols1314 = @("A","Q","R","Z","AB")
foreach ($col in $cols1314) {
    Swap-Cell 13 14 $col
}
